$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends at column AZ (52), which holds the "Mean"
# column (header + 13 data rows, r = 2..14). A new run ("Run 50") needs
# to be inserted as a data column right before "Mean", i.e. the old
# "Mean" column slides one column over to BA (53) and the new "Run 50"
# column takes its old place in AZ (52).

$oldMeanCol = 52   # AZ
$newRunCol  = 52   # AZ (after the shift, this holds "Run 50")
$newMeanCol = 53   # BA (after the shift, this holds "Mean")
$firstDataRow = 2
$lastDataRow  = 14

# New "Run 50" values (one per data row, r = 2..14)
$run50Values = @{
    2  = 15512923.52416456
    3  = 15512923.52416456
    4  = 15512923.52416456
    5  = 15512923.52416456
    6  = 15512923.52416456
    7  = 15512923.52416456
    8  = 15512923.52416456
    9  = 15512923.52416456
    10 = 15512923.52416456
    11 = 15512923.52416456
    12 = 15512923.52416456
    13 = 15512923.52416456
    14 = 15512923.52416456
}

# New recomputed "Mean" values for the extended row (r = 2..14)
$newMeanValues = @{
    2  = 86503685.70387085
    3  = 86503685.70387085
    4  = 86503685.70387085
    5  = 86503685.70387085
    6  = 86503685.70387085
    7  = 86503685.70387085
    8  = 86503685.70387085
    9  = 86503685.70387085
    10 = 86503685.70387085
    11 = 86503685.70387085
    12 = 86503685.70387085
    13 = 86503685.70387085
    14 = 86503685.70387085
}

# 1) Move the existing "Mean" header (with its bold/bordered style) from
#    AZ1 into the new BA1 cell first, while AZ1 still holds "Mean" -
#    Copy() duplicates both the value and the formatting.
$ws.Cells.Item(1, $oldMeanCol).Copy($ws.Cells.Item(1, $newMeanCol))

# 2) Now overwrite the header in AZ1 with the new "Run 50" label (it
#    already carries the same header style, so no extra formatting is
#    required).
$ws.Cells.Item(1, $newRunCol).Value = "Run 50"

# 3) Data rows: write the "Run 50" values into AZ, and the recomputed
#    "Mean" values into the new BA column.
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $ws.Cells.Item($r, $newRunCol).Value  = $run50Values[$r]
    $ws.Cells.Item($r, $newMeanCol).Value = $newMeanValues[$r]
}
